$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.458.95"
$ws.Range("E2").Value = "  +2.26%  "
$ws.Range("D3").Value = "3.389.89"
$ws.Range("E3").Value = "  +1.77%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +0.85%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +0.82%  "
$ws.Range("E9").Value = "  +5.67%  "
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("E12").Value = "  +2.96%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "681.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("D15").Value = "3.935.36"
$ws.Range("E15").Value = "  +1.78%  "
$ws.Range("D16").Value = "69.497.42"
$ws.Range("E16").Value = "  +2.30%  "
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("D18").Value = "3.389.87"
$ws.Range("E18").Value = "  +1.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.65"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.86%  "
$ws.Range("E21").Value = "  +0.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.67%  "
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("E26").Value = "  +1.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.71"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "34.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.55%  "
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "557.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.01%  "
$ws.Range("E33").Value = "  +4.92%  "
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.71%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").Value = "3.683.44"
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.28%  "
$ws.Range("E39").Value = "  +4.15%  "
$ws.Range("E40").Value = "  +3.43%  "
$ws.Range("E41").Value = "  +1.98%  "
$ws.Range("E42").Value = "  +2.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.340"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("E44").Value = "  +4.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.31"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.97%  "
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("E47").Value = "  +0.74%  "
$ws.Range("E48").Value = "  +5.20%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  +1.43%  "
$ws.Range("E51").Value = "  +2.94%  "
